$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the Thursday (quinta) Máquinas Elétricas entries at 8:40 and 9:50
$ws.Range("E4").Value = "-"
$ws.Range("E6").Value = "-"

# Move the class to Monday (segunda) at 10:40 and 11:30
$ws.Range("B7").Value = "MCT-3A-Máquinas Elétricas"
$ws.Range("B8").Value = "MCT-3A-Máquinas Elétricas"
